# Update crypto symbol data refreshed on 2023-02-11 by GitHub Actions.
# Rows 7 and 8 swap the MXToken / BTSEToken entries (ranking order changed),
# and D/E columns (Price / Volume(1h)) are refreshed with new values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.38"
$ws.Range("E2").Value = "'-0.28%"
$ws.Range("D3").Value = "'40.88"
$ws.Range("E3").Value = "'1.63%"
$ws.Range("D4").Value = "'5.107"
$ws.Range("E4").Value = "'-0.40%"
$ws.Range("D5").Value = "'0.07638"
$ws.Range("E5").Value = "'-1.52%"
$ws.Range("D6").Value = "'1.600"
$ws.Range("E6").Value = "'-1.73%"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "'0.9038"
$ws.Range("E7").Value = "'2.50%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.426"
$ws.Range("E8").Value = "'0.27%"
$ws.Range("E9").Value = "'8.28%"
$ws.Range("D10").Value = "'0.1778"
$ws.Range("E10").Value = "'1.69%"
$ws.Range("D11").Value = "'0.09149"
$ws.Range("E11").Value = "'1.23%"
$ws.Range("E12").Value = "'-5.11%"
$ws.Range("D13").Value = "'0.1051"
$ws.Range("E13").Value = "'-0.46%"
$ws.Range("D14").Value = "'0.001259"
$ws.Range("E14").Value = "'-0.17%"
$ws.Range("D15").Value = "'0.005671"
$ws.Range("E15").Value = "'-2.72%"
$ws.Range("D16").Value = "'3.350"
$ws.Range("E16").Value = "'-0.17%"
$ws.Range("D17").Value = "'4.244"
$ws.Range("E17").Value = "'-0.25%"
$ws.Range("E18").Value = "'0.51%"
$ws.Range("D19").Value = "'6.542"
$ws.Range("E19").Value = "'-6.90%"
$ws.Range("D20").Value = "'0.1364"
$ws.Range("E20").Value = "'1.87%"
$ws.Range("D21").Value = "'0.2829"
$ws.Range("E21").Value = "'1.51%"
$ws.Range("D22").Value = "'0.04064"
$ws.Range("E22").Value = "'-2.85%"
$ws.Range("D23").Value = "'0.001229"
$ws.Range("E23").Value = "'2.23%"
$ws.Range("D24").Value = "'0.004128"
$ws.Range("E24").Value = "'1.19%"
$ws.Range("E25").Value = "'0.03%"
$ws.Range("D26").Value = "'0.0003745"
$ws.Range("E26").Value = "'-95.00%"
$ws.Range("D38").Value = "'0.02414"
$ws.Range("E38").Value = "'1.58%"
$ws.Range("D39").Value = "'0.05179"
$ws.Range("E39").Value = "'-0.75%"
$ws.Range("D40").Value = "'0.007798"
$ws.Range("E40").Value = "'-1.65%"
$ws.Range("D41").Value = "'0.1304"
$ws.Range("E41").Value = "'-1.90%"
$ws.Range("D42").Value = "'0.007047"
$ws.Range("E42").Value = "'11.12%"
$ws.Range("E43").Value = "'-0.31%"
$ws.Range("D44").Value = "'0.008793"
$ws.Range("E44").Value = "'0.28%"
$ws.Range("D45").Value = "'0.3339"
$ws.Range("E45").Value = "'-0.66%"
$ws.Range("D46").Value = "'0.00006959"
$ws.Range("E46").Value = "'6.26%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.06%"
$ws.Range("D48").Value = "'0.03075"
$ws.Range("E48").Value = "'616.93%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.06%"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'-0.06%"
